$d = $word.ActiveDocument

$d.Content.Find.Execute("Projet Python | Trello", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://trello.com/b/ZtudX72t/projet-python", 2)

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*trello.com/b/ZtudX72t*") {
        $p.Range.Font.Reset()
    }
}
